$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1932629775376711
$ws.Range("C2").Value = 0.07381479721892248
$ws.Range("D2").Value = 0.02746420182479881
$ws.Range("F2").Value = 0.5352490448350977
$ws.Range("G2").Value = 0.3762415226960698
$ws.Range("H2").Value = 0.5434424720458182
$ws.Range("K2").Value = 0.1908499007080593
$ws.Range("M2").Value = 0.8869244610980616
$ws.Range("O2").Value = 1.780762920721443
$ws.Range("B3").Value = 0.1688055617385089
$ws.Range("C3").Value = 0.07314745209910711
$ws.Range("D3").Value = 0.02470966275431152
$ws.Range("F3").Value = 0.5345269578717549
$ws.Range("G3").Value = 0.3769269306848031
$ws.Range("H3").Value = 0.5471045313039937
$ws.Range("K3").Value = 0.1684584451220701
$ws.Range("M3").Value = 0.7896836777914586
$ws.Range("O3").Value = 1.789648995828045
$ws.Range("B4").Value = 0.1537374814116248
$ws.Range("C4").Value = 0.07274514631495776
$ws.Range("D4").Value = 0.0230057236765191
$ws.Range("F4").Value = 0.5344280057941688
$ws.Range("G4").Value = 0.3776248289870949
$ws.Range("H4").Value = 0.5495901237623286
$ws.Range("K4").Value = 0.1546261826340327
$ws.Range("M4").Value = 0.7305018596744048
$ws.Range("O4").Value = 1.796180741809621
$ws.Range("B5").Value = 0.1475848219285183
$ws.Range("C5").Value = 0.07258313150816775
$ws.Range("D5").Value = 0.02230822767126028
$ws.Range("F5").Value = 0.5344742284773076
$ws.Range("G5").Value = 0.3779787536163468
$ws.Range("H5").Value = 0.5506626492129953
$ws.Range("K5").Value = 0.1489689088776345
$ws.Range("M5").Value = 0.7065106605033549
$ws.Range("O5").Value = 1.799112682926491
$ws.Range("B6").Value = 0.1465624515200972
$ws.Range("C6").Value = 0.07255634716018022
$ws.Range("D6").Value = 0.0221922216483037
$ws.Range("F6").Value = 0.5344871292000022
$ws.Range("G6").Value = 0.3780417175342166
$ws.Range("H6").Value = 0.5508443429015699
$ws.Range("K6").Value = 0.1480283004237464
$ws.Range("M6").Value = 0.7025343692743178
$ws.Range("O6").Value = 1.799615840806155
$ws.Range("B7").Value = 0.1536545535272893
$ws.Range("C7").Value = 0.07274295344776505
$ws.Range("D7").Value = 0.02299632959454811
$ws.Range("F7").Value = 0.5344282788265602
$ws.Range("G7").Value = 0.3776293208234378
$ws.Range("H7").Value = 0.5496043467769383
$ws.Range("K7").Value = 0.1545499689795378
$ws.Range("M7").Value = 0.7301778041441338
$ws.Range("O7").Value = 1.796219189340974
$ws.Range("B8").Value = 0.1848410016210096
$ws.Range("C8").Value = 0.07358318396342867
$ws.Range("D8").Value = 0.02651709060397422
$ws.Range("F8").Value = 0.5349285492121538
$ws.Range("G8").Value = 0.3764202619071995
$ws.Range("H8").Value = 0.5446559526779993
$ws.Range("K8").Value = 0.183147053465845
$ws.Range("M8").Value = 0.8532833283828438
$ws.Range("O8").Value = 1.783603413812429
$ws.Range("B9").Value = 0.2455700092432949
$ws.Range("C9").Value = 0.07528762589881666
$ws.Range("D9").Value = 0.03331901840297746
$ws.Range("F9").Value = 0.5386458815275077
$ws.Range("G9").Value = 0.3762545264498272
$ws.Range("H9").Value = 0.5368328395680422
$ws.Range("K9").Value = 0.2385386174086079
$ws.Range("M9").Value = 1.099139928700623
$ws.Range("O9").Value = 1.767412900596142
$ws.Range("B10").Value = 0.2899018116995364
$ws.Range("C10").Value = 0.07657148864463892
$ws.Range("D10").Value = 0.03825182879519673
$ws.Range("F10").Value = 0.5430513374943615
$ws.Range("G10").Value = 0.3774874250627533
$ws.Range("H10").Value = 0.5322313919583905
$ws.Range("K10").Value = 0.2787897515464692
$ws.Range("M10").Value = 1.282900772118865
$ws.Range("O10").Value = 1.760750976294503
$ws.Range("B11").Value = 0.3100023649983257
$ws.Range("C11").Value = 0.07716179413920088
$ws.Range("D11").Value = 0.0404814259019588
$ws.Range("F11").Value = 0.545420482901001
$ws.Range("G11").Value = 0.3783447262719193
$ws.Range("H11").Value = 0.5303869223181295
$ws.Range("K11").Value = 0.2969991016747429
$ws.Range("M11").Value = 1.367271955679385
$ws.Range("O11").Value = 1.75886146233529
$ws.Range("B12").Value = 0.3176038839142734
$ws.Range("C12").Value = 0.07738617504170975
$ws.Range("D12").Value = 0.04132360312215155
$ws.Range("F12").Value = 0.5463702063068325
$ws.Range("G12").Value = 0.3787121652640764
$ws.Range("H12").Value = 0.5297242373030286
$ws.Range("K12").Value = 0.3038794476042312
$ws.Range("M12").Value = 1.399340452926154
$ws.Range("O12").Value = 1.758310397072904
$ws.Range("B13").Value = 0.3159672206930111
$ws.Range("C13").Value = 0.07733781391306138
$ws.Range("D13").Value = 0.04114232051702515
$ws.Range("F13").Value = 0.5461633265973447
$ws.Range("G13").Value = 0.3786311245759606
$ws.Range("H13").Value = 0.5298653671562903
$ws.Range("K13").Value = 0.3023983249259743
$ws.Range("M13").Value = 1.392428521435377
$ws.Range("O13").Value = 1.758421758776905
$ws.Range("B14").Value = 0.3106279528964535
$ws.Range("C14").Value = 0.0771802374779611
$ws.Range("D14").Value = 0.04055075520223284
$ws.Range("F14").Value = 0.54549756298087
$ws.Range("G14").Value = 0.3783740970033449
$ws.Range("H14").Value = 0.5303316857942946
$ws.Range("K14").Value = 0.2975654587448844
$ws.Range("M14").Value = 1.369907824402958
$ws.Range("O14").Value = 1.75881282783817
$ws.Range("B15").Value = 0.3073561608088085
$ws.Range("C15").Value = 0.07708382573541428
$ws.Range("D15").Value = 0.0401881262074113
$ws.Range("F15").Value = 0.5450966133524346
$ws.Range("G15").Value = 0.3782222387200989
$ws.Range("H15").Value = 0.530621978632567
$ws.Range("K15").Value = 0.2946031991588143
$ws.Range("M15").Value = 1.356128955179869
$ws.Range("O15").Value = 1.759073796659607
$ws.Range("B16").Value = 0.2885868093643467
$ws.Range("C16").Value = 0.07653303211515805
$ws.Range("D16").Value = 0.0381058254511828
$ws.Range("F16").Value = 0.5429038615780115
$ws.Range("G16").Value = 0.3774373781301961
$ws.Range("H16").Value = 0.532356938361346
$ws.Range("K16").Value = 0.2775976418800212
$ws.Range("M16").Value = 1.277403169529492
$ws.Range("O16").Value = 1.76089745282377
$ws.Range("B17").Value = 0.2770550329609307
$ws.Range("C17").Value = 0.07619670357248509
$ws.Range("D17").Value = 0.03682468181209941
$ws.Range("F17").Value = 0.5416522432027762
$ws.Range("G17").Value = 0.377031939158698
$ws.Range("H17").Value = 0.5334849949709195
$ws.Range("K17").Value = 0.2671389615720585
$ws.Range("M17").Value = 1.229311615470891
$ws.Range("O17").Value = 1.76230871716578
$ws.Range("B18").Value = 0.2704160645556897
$ws.Range("C18").Value = 0.07600384980369768
$ws.Range("D18").Value = 0.03608645339040351
$ws.Range("F18").Value = 0.5409667037533197
$ws.Range("G18").Value = 0.376826634512021
$ws.Range("H18").Value = 0.5341572366638871
$ws.Range("K18").Value = 0.2611139258652031
$ws.Range("M18").Value = 1.201723363944211
$ws.Range("O18").Value = 1.763227815665545
$ws.Range("B19").Value = 0.2681671782552826
$ws.Range("C19").Value = 0.07593865649545251
$ws.Range("D19").Value = 0.03583627207681417
$ws.Range("F19").Value = 0.5407404904771127
$ws.Range("G19").Value = 0.3767619074267898
$ws.Range("H19").Value = 0.5343888672596719
$ws.Range("K19").Value = 0.2590723425977899
$ws.Range("M19").Value = 1.192394733557251
$ws.Range("O19").Value = 1.763557437277044
$ws.Range("B20").Value = 0.278283255848379
$ws.Range("C20").Value = 0.07623244530491746
$ws.Range("D20").Value = 0.03696120171883166
$ws.Range("F20").Value = 0.5417819235716479
$ws.Range("G20").Value = 0.3770722107779676
$ws.Range("H20").Value = 0.533362488187386
$ws.Range("K20").Value = 0.2682532914001854
$ws.Range("M20").Value = 1.234423464027714
$ws.Range("O20").Value = 1.762147369899708
$ws.Range("B21").Value = 0.312196504852551
$ws.Range("C21").Value = 0.07722649901248246
$ws.Range("D21").Value = 0.04072457034529009
$ws.Range("F21").Value = 0.5456916862537824
$ws.Range("G21").Value = 0.3784484293704793
$ws.Range("H21").Value = 0.5301937457212489
$ws.Range("K21").Value = 0.2989854043761397
$ws.Range("M21").Value = 1.376519412466862
$ws.Range("O21").Value = 1.758693495184133
$ws.Range("B22").Value = 0.3343014816575476
$ws.Range("C22").Value = 0.07788106998851418
$ws.Range("D22").Value = 0.04317174920998923
$ws.Range("F22").Value = 0.5485534533365097
$ws.Range("G22").Value = 0.3795974037477663
$ws.Range("H22").Value = 0.5283313198058437
$ws.Range("K22").Value = 0.3189820921249691
$ws.Range("M22").Value = 1.470084181145126
$ws.Range("O22").Value = 1.757394919581174
$ws.Range("B23").Value = 0.3225092639202103
$ws.Range("C23").Value = 0.07753128345200366
$ws.Range("D23").Value = 0.0418667973918474
$ws.Range("F23").Value = 0.5469980017047362
$ws.Range("G23").Value = 0.3789612860947926
$ws.Range("H23").Value = 0.529306249105602
$ws.Range("K23").Value = 0.3083177822486789
$ws.Range("M23").Value = 1.420080726563626
$ws.Range("O23").Value = 1.758000145646719
$ws.Range("B24").Value = 0.2777280048802311
$ws.Range("C24").Value = 0.07621628488814736
$ws.Range("D24").Value = 0.03689948630303519
$ws.Range("F24").Value = 0.5417231890789083
$ws.Range("G24").Value = 0.3770539174395111
$ws.Range("H24").Value = 0.5334177996629208
$ws.Range("K24").Value = 0.267749540842928
$ws.Range("M24").Value = 1.232112210251969
$ws.Range("O24").Value = 1.762219979429091
$ws.Range("B25").Value = 0.2291895820961827
$ws.Range("C25").Value = 0.07482076384703618
$ws.Range("D25").Value = 0.03149009495993482
$ws.Range("F25").Value = 0.5373466630757164
$ws.Range("G25").Value = 0.3760622144887122
$ws.Range("H25").Value = 0.5387478942659598
$ws.Range("K25").Value = 0.2236300104039941
$ws.Range("M25").Value = 1.032110385368568
$ws.Range("O25").Value = 1.770875414199423
